$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "67.396.18"
$ws.Cells.Item(2, 5).Value = "  +2.02%  "
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "3.918.76"
$ws.Cells.Item(3, 5).Value = "  +3.65%  "
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "1.00"
$ws.Cells.Item(4, 5).Value = "  +0.13%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "469.92"
$ws.Cells.Item(5, 5).Value = "  +9.41%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "144.77"
$ws.Cells.Item(6, 5).Value = "  +4.60%  "
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.624"
$ws.Cells.Item(7, 5).Value = "  +0.44%  "
$ws.Cells.Item(8, 5).Value = "  -0.09%  "
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.732"
$ws.Cells.Item(9, 5).Value = "  -0.51%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.168"
$ws.Cells.Item(10, 5).Value = "  +11.29%  "
$ws.Cells.Item(11, 5).Value = "  +10.13%  "
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "43.31"
$ws.Cells.Item(12, 5).Value = "  +1.72%  "
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "4.548.20"
$ws.Cells.Item(13, 5).Value = "  +3.96%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "10.40"
$ws.Cells.Item(14, 5).Value = "  -0.20%  "
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "15.01"
$ws.Cells.Item(15, 5).Value = "  +0.06%  "
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "3.921.11"
$ws.Cells.Item(16, 5).Value = "  +3.66%  "
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "19.87"
$ws.Cells.Item(18, 5).Value = "  -0.39%  "
$ws.Cells.Item(19, 5).Value = "  +3.06%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "67.657.60"
$ws.Cells.Item(20, 5).Value = "  +2.31%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "433.80"
$ws.Cells.Item(21, 5).Value = "  +7.27%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "14.63"
$ws.Cells.Item(22, 5).Value = "  -2.04%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "3.34"
$ws.Cells.Item(23, 5).Value = "  +1.85%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "87.84"
$ws.Cells.Item(24, 5).Value = "  +3.85%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "38.65"
$ws.Cells.Item(25, 5).Value = "  +5.73%  "
$ws.Cells.Item(26, 5).Value = "  +7.24%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "5.75"
$ws.Cells.Item(27, 5).Value = "  +4.24%  "
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "10.21"
$ws.Cells.Item(28, 5).Value = "  +3.66%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "9.62"
$ws.Cells.Item(29, 5).Value = "  -4.02%  "
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "723.27"
$ws.Cells.Item(30, 5).Value = "  +2.80%  "
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "13.61"
$ws.Cells.Item(31, 5).Value = "  -0.97%  "
$ws.Cells.Item(32, 5).Value = "  -2.96%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "2.83"
$ws.Cells.Item(33, 5).Value = "  +4.72%  "
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "42.86"
$ws.Cells.Item(34, 5).Value = "  +4.14%  "
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "0.154"
$ws.Cells.Item(35, 5).Value = "  +3.68%  "
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "57.79"
$ws.Cells.Item(36, 5).Value = "  +2.92%  "
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.0₃0799"
$ws.Cells.Item(37, 5).Value = "  +19.01%  "
$ws.Cells.Item(38, 5).Value = "  -0.11%  "
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "5.37"
$ws.Cells.Item(39, 5).Value = "  -7.85%  "
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.0476"
$ws.Cells.Item(40, 5).Value = "  +1.11%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "3.07"
$ws.Cells.Item(41, 5).Value = "  +1.79%  "
$ws.Cells.Item(42, 2).Value = "Fetch.AI"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "2.58"
$ws.Cells.Item(42, 5).Value = "  -8.34%  "
$ws.Cells.Item(43, 2).Value = "Stellar"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.141"
$ws.Cells.Item(43, 5).Value = "  +0.16%  "
$ws.Cells.Item(44, 2).Value = "FirstDigitalUSD"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "1.00"
$ws.Cells.Item(44, 5).Value = "  +0.09%  "
$ws.Cells.Item(45, 2).Value = "TheGraph"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.336"
$ws.Cells.Item(45, 5).Value = "  +1.93%  "
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "2.79"
$ws.Cells.Item(46, 5).Value = "  +4.32%  "
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "2.17"
$ws.Cells.Item(47, 5).Value = "  +5.47%  "
$ws.Cells.Item(48, 5).Value = "  +1.50%  "
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "3.15"
$ws.Cells.Item(49, 5).Value = "  -0.96%  "
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "145.29"
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "2.88"
$ws.Cells.Item(51, 5).Value = "  +3.34%  "
